# Changed the logo and added a new job.
# Replace the "Home Job2 (INDIA) PVT. LTD." / Gurugaon row (row 3) with a new
# job posting for "Itel MNC Company" in Noida Sector 63, and drop the now
# unused trailing detail columns (O3:S3) that the new posting doesn't use.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = "Itel MNC Company"
$ws.Range("C3").Value = "Noida Sector 63"
$ws.Range("D3").Value = "Not in use"
$ws.Range("E3").Value = "Not in use"
$ws.Range("F3").Value = "इंटरव्यू पता :- ARK Workforce Plot Number – 246 Block – G Sector-63 Noida"
# H3 (time) is assigned before G3 (date) so the new shared-string table keeps
# the same entry order as the source workbook.
$ws.Range("H3").Value = "इंटरव्यू टाइम :- 11 बजे तक "
$ws.Range("G3").Value = "इंटरव्यू की तारीख: 04/04/2025"
$ws.Range("I3").Value = "Contact person: Ankit -8477873797, Vishan – 9315473717"
$ws.Range("J3").Value = "Qualification:-10th 12th graduate ITI or Diploma"
$ws.Range("K3").Value = "Salary: 11000 Rs in hand"
$ws.Range("L3").Value = "Total Vacancy :- 100(Only Boys)"
$ws.Range("M3").Value = "अतिरिक्त लाभ: attendance award: 1000 Rs, Night Alowance 50  Rs Per Night, Overtime: 100 Rs Per Hour, Lunch_Dinner_Free"
$ws.Range("N3").Value = "NOTE-: एक साल के बाद भारत सरकार के द्वारा APPPRENTICESHIP का CERTIFICATE मिलता है जिसको कही भी आप एक्सपेरिएंस सर्टिफिकेट की तरह उपयोग कर सकते है."

# The previous job posting used columns up to S3; the new posting only uses
# up to N3, so remove the now-stale trailing cells.
$ws.Range("O3:S3").ClearContents()

# The view had scrolled and a different cell became the active selection.
$null = $ws.Range("U3").Select()
